$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: "Esta es una frase de prueba creada por Kevin." becomes three
# runs: "Esta es una frase de prueba creada por Kevin" + ", y editada" + "."
#
# A plain InsertAfter/InsertBefore on a same-formatted Range gets silently
# coalesced back into a single run, so the insertion is done with
# TrackRevisions on (each tracked insertion keeps its own run) and then
# accepted, which leaves the run boundaries in place without leaving any
# revision markup or formatting residue behind.
# -----------------------------------------------------------------------
$d.TrackRevisions = $true

$found = $d.Content
$found.Find.Execute("Kevin.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$period = $d.Range($found.End - 1, $found.End)
$period.InsertBefore(", y editada")

$d.TrackRevisions = $false
$d.AcceptAllRevisions()

# -----------------------------------------------------------------------
# Part 2: add a new paragraph right after that sentence:
# "Esta es otra frase para probar el commit"
# -----------------------------------------------------------------------
$sentence = $d.Content
$sentence.Find.Execute("Esta es una frase de prueba creada por Kevin, y editada.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the paragraph index for the sentence we just found (re-reading
# from $d.Paragraphs keeps this accurate even if Range bookkeeping on
# $sentence itself goes stale after further edits).
$sentenceParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i).Range
    if ($p.Start -le $sentence.Start -and $sentence.Start -lt $p.End) {
        $sentenceParaIndex = $i
    }
}

$sentence.InsertParagraphAfter()

$newPara = $d.Paragraphs($sentenceParaIndex + 1).Range
$newPara.Text = "Esta es otra frase para probar el commit"
